# Apply "Trade #30 closed" update to the live trading results workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet - update aggregate stats
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.38   # Current Capital
$summary.Range("B4").Value = -0.62     # Total P&L $
$summary.Range("B6").Value = 30        # Total Trades
$summary.Range("B8").Value = 15        # Losing Trades
$summary.Range("B9").Value = 26.67     # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - update MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.38      # Capital
$status.Range("D4").Value = 30         # Trades
$status.Range("E4").Value = -0.62      # P&L $
$status.Range("F4").Value = -0.62      # P&L %
$status.Range("G4").Value = 26.67      # Win Rate %

# ---------------------------------------------------------------------
# 3) New trade row (Trade #30) appended to "All Trades" and
#    "MarketMaking" sheets (row 31 on both).
# ---------------------------------------------------------------------
function Set-TradeRow($ws) {
    $ws.Range("A31").Value = 30
    # Leading apostrophe forces these to stay plain text (matching the other
    # rows) instead of being auto-parsed into date/time serial values.
    $ws.Range("B31").Value = "'2026-02-17"
    $ws.Range("C31").Value = "'08:03:45"
    $ws.Range("D31").Value = "MarketMaking"
    $ws.Range("E31").Value = "UP"
    $ws.Range("F31").Value = 0.3
    $ws.Range("G31").Value = 0.277228
    $ws.Range("H31").Value = "CLOSED"
    $ws.Range("I31").Value = -7.5908
    $ws.Range("J31").Value = -0.02
    $ws.Range("K31").Value = 99.38
    $ws.Range("L31").Value = 0
    $ws.Range("M31").Value = 0
    $ws.Range("N31").Value = 0.6
    $ws.Range("O31").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P31").Value = "early_exit"
    $ws.Range("Q31").Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Set-TradeRow $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Set-TradeRow $marketMaking
